$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.14%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-6.40%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.060"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.54%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07651"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-6.46%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.247"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.60%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.602"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-9.25%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9146"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.30%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1039"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-7.61%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1780"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.96%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09423"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.64%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04441"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.00%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.02%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001263"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.05%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005821"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.06%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2,406.21%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.361"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.32%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.419"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.47%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3317"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.944"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-7.37%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1349"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.20%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2817"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "10.17%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04157"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.61%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001207"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.50%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004107"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.04%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001304"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "6.56%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02458"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-6.80%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05175"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-7.80%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007910"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.09%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1316"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-5.98%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007091"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.26%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001954"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.69%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007419"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.43%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3066"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-11.96%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006446"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-5.16%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.23%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003004"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-26.99%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004548"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "35.59%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.23%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.23%"
